# Auto-generated PowerShell Excel COM-interop script
# Applies metadata field updates and appends new rows (26-29) to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$fileName = "KRTINN-Openscienceenonderzoeksinformatie-060923-1508-1296.pdf"

$ws.Cells.Item(2, 1).Value = $fileName
$ws.Cells.Item(2, 2).Value = "- Title"
$ws.Cells.Item(2, 3).Value = "Open science en onderzoeksinformatie"

$ws.Cells.Item(3, 1).Value = $fileName
$ws.Cells.Item(3, 2).Value = "- Creator (Author)"
$ws.Cells.Item(3, 3).Value = "John Doove, Germaine Poot, Karin van Grieken, Gül Akcaova"

$ws.Cells.Item(4, 1).Value = $fileName
$ws.Cells.Item(4, 2).Value = "- Date"
$ws.Cells.Item(4, 3).Value = "No specific date mentioned"

$ws.Cells.Item(5, 1).Value = $fileName
$ws.Cells.Item(5, 2).Value = "- Description"
$ws.Cells.Item(5, 3).Value = "The text describes the ambitions and initiatives related to open science and research information in the Netherlands, particularly in the context of higher education institutions (hbo). It mentions projects, platforms, and collaborations aimed at making research information open, accessible, and visible. It also discusses the role of SURF (a collaborative organization for ICT in Dutch education and research) in facilitating open science in hbo."

$ws.Cells.Item(6, 1).Value = $fileName
$ws.Cells.Item(6, 2).Value = "- Keywords"
$ws.Cells.Item(6, 3).Value = "Open science, onderzoeksinformatie, ambitie, hogescholen, SURF, project, NPPO, Knowledge Exchange, Persistent Identifiers, ORCID, roadmap"

$ws.Cells.Item(7, 1).Value = $fileName
$ws.Cells.Item(7, 2).Value = "- Type"
$ws.Cells.Item(7, 3).Value = "Text (document, article)"

$ws.Cells.Item(8, 1).Value = $fileName
$ws.Cells.Item(8, 2).Value = "- Identifier (DOI)"
$ws.Cells.Item(8, 3).Value = "No DOI provided"

$ws.Cells.Item(9, 1).Value = $fileName
$ws.Cells.Item(9, 2).Value = "- Publisher"
$ws.Cells.Item(9, 3).Value = "No specific publisher mentioned"

$ws.Cells.Item(10, 1).Value = $fileName
$ws.Cells.Item(10, 2).Value = "- Rights"
$ws.Cells.Item(10, 3).Value = "No specific rights mentioned"

$ws.Cells.Item(11, 1).Value = $fileName
$ws.Cells.Item(11, 2).Value = "- Language"
$ws.Cells.Item(11, 3).Value = "Dutch"

$ws.Cells.Item(12, 1).Value = $fileName
$ws.Cells.Item(12, 2).Value = "- Format"
$ws.Cells.Item(12, 3).Value = "Textual format (possibly in a digital document)"

$ws.Cells.Item(13, 1).Value = $fileName
$ws.Cells.Item(13, 2).Value = "- Source"
$ws.Cells.Item(13, 3).Value = "No specific source mentioned"

$ws.Cells.Item(14, 1).Value = $fileName
$ws.Cells.Item(14, 2).Value = "- Relation"
$ws.Cells.Item(14, 3).Value = "Various projects and collaborations are mentioned, such as the relation between the Open Science and onderzoeksinformatie project and the NPPO project, and the relation between SURF and the `"Adviescollege Open Science in het hbo`" of the Vereniging Hogescholen."

$ws.Cells.Item(15, 1).Value = $fileName
$ws.Cells.Item(15, 2).Value = "- Coverage"
$ws.Cells.Item(15, 3).Value = "The text primarily focuses on open science and research information initiatives in the Netherlands, particularly in hbo (higher education). It also mentions European collaborations through Knowledge Exchange."

$ws.Cells.Item(16, 1).Value = $fileName
$ws.Cells.Item(16, 2).Value = "- Title"
$ws.Cells.Item(16, 3).Value = "Open science en onderzoeksinformatie"

$ws.Cells.Item(17, 1).Value = $fileName
$ws.Cells.Item(17, 2).Value = "- Creator (Author)"
$ws.Cells.Item(17, 3).Value = "John Doove, Germaine Poot, Karin van Grieken, Gül Akcaova"

$ws.Cells.Item(18, 1).Value = $fileName
$ws.Cells.Item(18, 2).Value = "- Date"
$ws.Cells.Item(18, 3).Value = "No specific date mentioned"

$ws.Cells.Item(19, 1).Value = $fileName
$ws.Cells.Item(19, 2).Value = "- Description"
$ws.Cells.Item(19, 3).Value = "The text describes the ambitions and initiatives related to open science and research information in the Netherlands, particularly in the context of higher education institutions (hbo). It mentions projects, platforms, and collaborations aimed at making research information open, accessible, and visible. It also discusses the role of SURF (a collaborative organization for ICT in Dutch education and research) in facilitating open science in hbo."

$ws.Cells.Item(20, 1).Value = $fileName
$ws.Cells.Item(20, 2).Value = "- Keywords"
$ws.Cells.Item(20, 3).Value = "Open science, onderzoeksinformatie, ambitie, hogescholen, SURF, project, NPPO, Knowledge Exchange, Persistent Identifiers, ORCID, roadmap"

$ws.Cells.Item(21, 1).Value = $fileName
$ws.Cells.Item(21, 2).Value = "- Type"
$ws.Cells.Item(21, 3).Value = "Text (document, article)"

$ws.Cells.Item(22, 1).Value = $fileName
$ws.Cells.Item(22, 2).Value = "- Identifier (DOI)"
$ws.Cells.Item(22, 3).Value = "No DOI provided"

$ws.Cells.Item(23, 1).Value = $fileName
$ws.Cells.Item(23, 2).Value = "- Publisher"
$ws.Cells.Item(23, 3).Value = "No specific publisher mentioned"

$ws.Cells.Item(24, 1).Value = $fileName
$ws.Cells.Item(24, 2).Value = "- Rights"
$ws.Cells.Item(24, 3).Value = "No specific rights mentioned"

$ws.Cells.Item(25, 1).Value = $fileName
$ws.Cells.Item(25, 2).Value = "- Language"
$ws.Cells.Item(25, 3).Value = "Dutch"

$ws.Cells.Item(26, 1).Value = $fileName
$ws.Cells.Item(26, 2).Value = "- Format"
$ws.Cells.Item(26, 3).Value = "Textual format (possibly in a digital document)"

$ws.Cells.Item(27, 1).Value = $fileName
$ws.Cells.Item(27, 2).Value = "- Source"
$ws.Cells.Item(27, 3).Value = "No specific source mentioned"

$ws.Cells.Item(28, 1).Value = $fileName
$ws.Cells.Item(28, 2).Value = "- Relation"
$ws.Cells.Item(28, 3).Value = "Various projects and collaborations are mentioned, such as the relation between the Open Science and onderzoeksinformatie project and the NPPO project, and the relation between SURF and the `"Adviescollege Open Science in het hbo`" of the Vereniging Hogescholen."

$ws.Cells.Item(29, 1).Value = $fileName
$ws.Cells.Item(29, 2).Value = "- Coverage"
$ws.Cells.Item(29, 3).Value = "The text primarily focuses on open science and research information initiatives in the Netherlands, particularly in hbo (higher education). It also mentions European collaborations through Knowledge Exchange."
